$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "25.768.73" or "1.002" that Excel
# would otherwise auto-parse as numbers/dates; force Text format, assign,
# then drop back to the default style so no stray "s" attribute is left
# on the cell (matches the un-styled inline-string cells in the source file).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.768.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.729.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5156'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2746'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.24'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06093'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.734.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6310'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.484'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.800.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006607'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.955.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.137'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('E24').Value = '  +4.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.114'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.503'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.756'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '101.73'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08271'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.659'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.455'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04468'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.617'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9685'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6103'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.654'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01569'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.922'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3795'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.984'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7145'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05359'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.167'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.548'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.88%  '
